$d = $word.ActiveDocument

# Table 1: Median (Q1-Q3) age values -> Mean (SE) values
# "54 (47 – 60)" appears twice and both map to "52 (1)" -> replace all
$d.Content.Find.Execute("54 (47 – 60)", $true, $false, $false, $false, $false, $true, 1, $false, "52 (1)", 2) | Out-Null

$d.Content.Find.Execute("54 (48 – 59)", $true, $false, $false, $false, $false, $true, 1, $false, "52 (1)", 2) | Out-Null

$d.Content.Find.Execute("52 (43 – 58)", $true, $false, $false, $false, $false, $true, 1, $false, "50 (2)", 2) | Out-Null

# Footer text for first table
$d.Content.Find.Execute("n (unweighted) (%); Median (Q1 – Q3)", $true, $false, $false, $false, $false, $true, 1, $false, "n (unweighted) (%); Mean (SE)", 2) | Out-Null

# Table 2 (second statistics table)
$d.Content.Find.Execute("123 (110 – 134)", $true, $false, $false, $false, $false, $true, 1, $false, "124 (2)", 2) | Out-Null

$d.Content.Find.Execute("133 (116 – 144)", $true, $false, $false, $false, $false, $true, 1, $false, "133 (2)", 2) | Out-Null

$d.Content.Find.Execute("124 (112 – 139)", $true, $false, $false, $false, $false, $true, 1, $false, "128 (3)", 2) | Out-Null

$d.Content.Find.Execute("124 (111 – 137)", $true, $false, $false, $false, $false, $true, 1, $false, "127 (1)", 2) | Out-Null

# Table 3 (third statistics table)
# "73 (67 – 81)" appears twice and both map to "74 (1)" -> replace all
$d.Content.Find.Execute("73 (67 – 81)", $true, $false, $false, $false, $false, $true, 1, $false, "74 (1)", 2) | Out-Null

$d.Content.Find.Execute("75 (67 – 85)", $true, $false, $false, $false, $false, $true, 1, $false, "76 (1)", 2) | Out-Null

$d.Content.Find.Execute("75 (66 – 82)", $true, $false, $false, $false, $false, $true, 1, $false, "75 (2)", 2) | Out-Null

# Footer text for second table
$d.Content.Find.Execute("n (unweighted)/N (unweighted) (%); Median (Q1 – Q3)", $true, $false, $false, $false, $false, $true, 1, $false, "n (unweighted)/N (unweighted) (%); Mean (SE)", 2) | Out-Null
